$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the four long student names with their short local-version names.
$ws.Range("B37").Value = "congying"
$ws.Range("B38").Value = "irfan"
$ws.Range("B39").Value = "jessica"
$ws.Range("B40").Value = "vic"

# Update the view: scroll position and selection, matching the saved view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B41").Select()

# Resize / reposition the workbook window (maximized-ish geometry from the diff).
$excel.ActiveWindow.WindowState = -4137
$excel.ActiveWindow.Top = -110
$excel.ActiveWindow.Left = -110
$excel.ActiveWindow.Width = 19420
$excel.ActiveWindow.Height = 10300
